$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column C (shifts Order/Active/Category/Image one column right)
$ws.Range("C1").EntireColumn.Insert()

# Set the header for the newly inserted column
$ws.Range("C1").Value = "SpecialPrice"

# Update selection to match the authored diff (C1 active cell)
$ws.Range("C1").Select()
